# Update the "Greeting" value for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE",
# matching the jgit commit that updated the Main.xlsx test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection recorded in the sheet view.
$ws.Activate()
$ws.Range("E8").Select()
